$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint")

# Update status dropdowns for tasks in rows 9 and 10 from "To do" to "Done"
$ws.Range("F9").Value = "Done"
$ws.Range("F10").Value = "Done"

# Update "Day 2" effort entries for rows 9 and 10
$ws.Range("T9").Value = 5
$ws.Range("T10").Value = 2

# Update selection to match the author's final cursor position
$ws.Range("T11").Select()
